# Table24.xlsx — fill in the "Region of residence" (col A) and "Nativity"
# (col B) labels that were left blank in several rows, per the source data
# wrangling pass ("progress in wrangling nfishers data").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel alignment constants used below.
$xlLeft    = -4131
$xlTop     = -4160
$xlJustify = -4130

function Set-LeftTop($cell) {
    $cell.HorizontalAlignment = $xlLeft
    $cell.VerticalAlignment = $xlTop
}

function Set-Justify($cell) {
    $cell.HorizontalAlignment = $xlJustify
}

function Set-WrapIndent($cell) {
    $cell.HorizontalAlignment = $xlLeft
    $cell.WrapText = $true
    $cell.IndentLevel = 1
}

# ---- "Eureka" block (rows 2-4): style stays left/top ----
Set-LeftTop $ws.Cells.Item(2, 1)
$ws.Cells.Item(2, 1).Value = "Eureka"

$ws.Cells.Item(3, 1).Value = "Eureka"

$ws.Rows.Item(4).RowHeight = 17
$ws.Cells.Item(4, 1).Value = "Eureka"

# ---- "Sacramento" block (rows 6-8) ----
$ws.Cells.Item(6, 1).Value = "Sacramento"
$ws.Cells.Item(7, 1).Value = "Sacramento"
$ws.Cells.Item(8, 1).Value = "Sacramento"

# ---- "San Francisco" block (rows 10-13) ----
$ws.Rows.Item(10).RowHeight = 17
$ws.Cells.Item(10, 1).Value = "San Francisco"
$ws.Cells.Item(11, 1).Value = "San Francisco"
$ws.Cells.Item(12, 1).Value = "San Francisco"
$ws.Cells.Item(13, 1).Value = "San Francisco"
$ws.Cells.Item(13, 2).Value = "Jugoslavia"

# ---- "Monterey" block (rows 15-18) ----
$ws.Rows.Item(15).RowHeight = 17
$ws.Cells.Item(15, 1).Value = "Monterey"
$ws.Cells.Item(16, 1).Value = "Monterey"
$ws.Cells.Item(16, 2).Value = "Italy"
$ws.Cells.Item(17, 1).Value = "Monterey"
$ws.Cells.Item(18, 1).Value = "Monterey"
$ws.Cells.Item(18, 2).Value = "Portugal"

# ---- "Santa Barbara" block (rows 20-22): justify alignment ----
$a20 = $ws.Cells.Item(20, 1)
Set-Justify $a20
$a20.Value = "Santa Barbara "

$ws.Rows.Item(21).RowHeight = 17
$a21 = $ws.Cells.Item(21, 1)
Set-Justify $a21
$a21.Value = "Santa Barbara "
$ws.Cells.Item(21, 2).Value = "Philippine Island"

$ws.Rows.Item(22).RowHeight = 17
$a22 = $ws.Cells.Item(22, 1)
Set-Justify $a22
$a22.Value = "Santa Barbara "
$ws.Cells.Item(22, 2).Value = "Portugal"

# ---- "Los Angeles" block (rows 24-31): justify alignment ----
$a24 = $ws.Cells.Item(24, 1)
Set-Justify $a24
$a24.Value = "Los Angeles"

$ws.Cells.Item(25, 1).Value = "Los Angeles"
$ws.Cells.Item(26, 1).Value = "Los Angeles"

$ws.Rows.Item(27).RowHeight = 17
$ws.Cells.Item(27, 1).Value = "Los Angeles"

$ws.Rows.Item(28).RowHeight = 17
$a28 = $ws.Cells.Item(28, 1)
Set-Justify $a28
$a28.Value = "Los Angeles"

$ws.Rows.Item(29).RowHeight = 17
$ws.Cells.Item(29, 1).Value = "Los Angeles"

$ws.Cells.Item(30, 1).Value = "Los Angeles"
$ws.Cells.Item(31, 1).Value = "Los Angeles"

# ---- "San Diego" block (rows 33-37): justify alignment ----
$a33 = $ws.Cells.Item(33, 1)
Set-Justify $a33
$a33.Value = "San Diego"

$ws.Rows.Item(34).RowHeight = 17
$a34 = $ws.Cells.Item(34, 1)
Set-Justify $a34
$a34.Value = "San Diego"
$ws.Cells.Item(34, 2).Value = "Portugal"

$ws.Cells.Item(35, 1).Value = "San Diego"

$ws.Rows.Item(36).RowHeight = 17
$ws.Cells.Item(36, 1).Value = "San Diego"

$ws.Cells.Item(37, 1).Value = "San Diego"

# ---- "Alaska Oregon Washington ..." block (rows 39-41): wrap + indent ----
$ws.Rows.Item(39).RowHeight = 51
$a39 = $ws.Cells.Item(39, 1)
Set-WrapIndent $a39
$a39.Value = "Alaska Oregon Washington and other states licensed in California"

$ws.Rows.Item(40).RowHeight = 51
$a40 = $ws.Cells.Item(40, 1)
Set-WrapIndent $a40
$a40.Value = "Alaska Oregon Washington and other states licensed in California"
$ws.Cells.Item(40, 2).Value = "Norway"

$ws.Rows.Item(41).RowHeight = 51
$a41 = $ws.Cells.Item(41, 1)
Set-WrapIndent $a41
$a41.Value = "Alaska Oregon Washington and other states licensed in California"

# ---- viewport: scrolled down, last click landed on A42 ----
$ws.Range("A42").Select()
